$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.806894511460617
$ws.Range("H2").Value = 0.0902626890121446
$ws.Range("K2").Value = -1.91805909292759

$ws.Range("G3").Value = 0.645651656549176
$ws.Range("H3").Value = 0.00917148459720282
$ws.Range("K3").Value = -30.7986636860888
$ws.Range("L3").Value = 0

$ws.Range("G4").Value = 1.04661626773175
$ws.Range("H4").Value = 0.0264884636899569
$ws.Range("K4").Value = 1.800266934792

$ws.Range("G5").Value = 0.800168606154522
$ws.Range("H5").Value = 0.0897055388561325
$ws.Range("K5").Value = -1.98854878730234

$ws.Range("G6").Value = 0.77095544598149
$ws.Range("H6").Value = 0.0878427688068737
$ws.Range("K6").Value = -2.28299440588269

$ws.Range("G7").Value = 0.616894344618247
$ws.Range("H7").Value = 0.0165133081613456
$ws.Range("K7").Value = -18.0457751574743
$ws.Range("L7").Value = 0

$ws.Range("G8").Value = 0.987502576858034
$ws.Range("H8").Value = 0.175265198353641
$ws.Range("K8").Value = -0.0708583512201663

$ws.Range("G9").Value = 0.901742026018586
$ws.Range("H9").Value = 0.115335475603683
$ws.Range("K9").Value = -0.808634927912254

$ws.Range("G11").Value = 0.913154099189985
$ws.Range("H11").Value = 0.135922891322433
$ws.Range("K11").Value = -0.610350646178583

$ws.Range("G14").Value = 1.65150691072196
$ws.Range("H14").Value = 0.361630369345433
$ws.Range("K14").Value = 2.29112795154251

$ws.Range("G15").Value = 1.15057478759627
$ws.Range("H15").Value = 0.205530046186997
$ws.Range("K15").Value = 0.785196623911994

$ws.Range("G17").Value = 0.696681788084856
$ws.Range("H17").Value = 0.0921145610371027
$ws.Range("K17").Value = -2.73354472590943
$ws.Range("L17").Value = 0.959

$ws.Range("G20").Value = 0.630945469873289
$ws.Range("H20").Value = 0.0182531305009342
$ws.Range("K20").Value = -15.9190776160157
$ws.Range("L20").Value = 0

$ws.Range("G21").Value = 0.602582653033566
$ws.Range("H21").Value = 0.0153422796469165
$ws.Range("K21").Value = -19.8944656948968
$ws.Range("L21").Value = 0

$ws.Range("G22").Value = 1.20678002705247
$ws.Range("H22").Value = 0.0252774723515188
$ws.Range("K22").Value = 8.97325313305902
$ws.Range("L22").Value = 0

$ws.Range("G23").Value = 0.955047118659211
$ws.Range("H23").Value = 0.0258645772188469
$ws.Range("K23").Value = -1.6983463755943

$ws.Range("G24").Value = 0.52283386841789
$ws.Range("H24").Value = 0.0141458938158284
$ws.Range("K24").Value = -23.9683213141566
$ws.Range("L24").Value = 0

$ws.Range("G25").Value = 0.499330979569955
$ws.Range("H25").Value = 0.0109918838313895
$ws.Range("K25").Value = -31.5485897246112
$ws.Range("L25").Value = 0

$ws.Range("G26").Value = 0.793907242088408
$ws.Range("H26").Value = 0.0832131928513175
$ws.Range("K26").Value = -2.20187175687539

$ws.Range("G27").Value = 0.765597779505446
$ws.Range("H27").Value = 0.0121167974437633
$ws.Range("K27").Value = -16.8765629973555
$ws.Range("L27").Value = 0

$ws.Range("G28").Value = 1.23597854434478
$ws.Range("H28").Value = 0.0190179257093596
$ws.Range("K28").Value = 13.7690159415274
$ws.Range("L28").Value = 0

$ws.Range("G29").Value = 0.964341599267324
$ws.Range("H29").Value = 0.100717242051771
$ws.Range("K29").Value = -0.347655920942992

$ws.Range("G30").Value = 0.64233092533922
$ws.Range("H30").Value = 0.0670546924226365
$ws.Range("K30").Value = -4.24025273150399
$ws.Range("L30").Value = 0.003

$ws.Range("G31").Value = 0.619426431800484
$ws.Range("H31").Value = 0.00835245121817608
$ws.Range("K31").Value = -35.5202689021657
$ws.Range("L31").Value = 0
